$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2023-07-16)
$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.000002220651329265522
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 23.48774463709945

# Row 3 (2023-03-10)
$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 2.652525447291612
